# Normalize the "Recorded By" column (column G) on the active sheet:
# for each comma-separated list of recorders, move any exact-case "System"
# entries to the end of the list and alphabetically (case-insensitively)
# sort the remaining entries.
#
# Example: "System, backup@backdoor.com, system" -> "backup@backdoor.com, system, System"
#          "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#          "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ",\s*"

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    $sortedOthers = @($otherParts | Sort-Object { $_.ToLower() })

    $newParts = @()
    $newParts += $sortedOthers
    $newParts += $systemParts
    $newVal = [string]::Join(", ", $newParts)

    if (-not $newVal.Equals($val)) {
        $cell.Value2 = $newVal
    }
}
